$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Value)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

# Cells whose new text looks like a plain number (single decimal point).
# Force them to remain text so Excel does not coerce them into numeric values.
Set-TextValue $ws 'D5' '0.7172'
Set-TextValue $ws 'D6' '239.11'
Set-TextValue $ws 'D8' '0.07831'
Set-TextValue $ws 'D9' '0.3074'
Set-TextValue $ws 'D10' '25.30'
Set-TextValue $ws 'D11' '0.08251'
Set-TextValue $ws 'D13' '0.7245'
Set-TextValue $ws 'D14' '5.247'
Set-TextValue $ws 'D15' '90.71'
Set-TextValue $ws 'D17' '5.861'
Set-TextValue $ws 'D18' '0.000007875'
Set-TextValue $ws 'D19' '242.48'
Set-TextValue $ws 'D20' '13.30'
Set-TextValue $ws 'D24' '7.772'
Set-TextValue $ws 'D25' '0.1551'
Set-TextValue $ws 'D26' '163.30'
Set-TextValue $ws 'D27' '9.008'
Set-TextValue $ws 'D28' '18.35'
Set-TextValue $ws 'D29' '1.934'
Set-TextValue $ws 'D30' '1.357'
Set-TextValue $ws 'D31' '1.483'
Set-TextValue $ws 'D33' '4.093'
Set-TextValue $ws 'D34' '0.05252'
Set-TextValue $ws 'D36' '0.7178'
Set-TextValue $ws 'D37' '1.003'
Set-TextValue $ws 'D38' '2.680'
Set-TextValue $ws 'D39' '0.01864'
Set-TextValue $ws 'D40' '2.725'
Set-TextValue $ws 'D42' '0.9088'
Set-TextValue $ws 'D43' '6.020'
Set-TextValue $ws 'D44' '72.34'
Set-TextValue $ws 'D45' '0.4310'
Set-TextValue $ws 'D47' '102.39'
Set-TextValue $ws 'D48' '0.5363'
Set-TextValue $ws 'D50' '9.159'
Set-TextValue $ws 'D51' '7.024'

# Cells that are safe to assign directly (percentages, links, names, multi-dot price strings).
$ws.Range('D2').Value = '29.414.78'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '1.873.71'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +0.79%  '
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -3.56%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('E10').Value = '  +9.04%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').Value = '1.884.40'
$ws.Range('E12').Value = '  +1.66%  '
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('E15').Value = '  +1.45%  '
$ws.Range('D16').Value = '29.448.55'
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('E19').Value = '  +2.11%  '
$ws.Range('E20').Value = '  -0.37%  '
$ws.Range('D21').Value = '2.125.32'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +4.86%  '
$ws.Range('E25').Value = '  +6.12%  '
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('E27').Value = '  +0.50%  '
$ws.Range('E28').Value = '  +1.50%  '
$ws.Range('E29').Value = '  -1.35%  '
$ws.Range('E30').Value = '  -4.95%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('E35').Value = '  +2.82%  '
$ws.Range('E36').Value = '  +1.37%  '
$ws.Range('E37').Value = '  +0.42%  '
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('D41').Value = '1.181.08'
$ws.Range('E41').Value = '  +3.41%  '
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E43').Value = '  +2.59%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E44').Value = '  +3.29%  '
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('E46').Value = '  +0.13%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('E50').Value = '  -0.37%  '
$ws.Range('E51').Value = '  +1.00%  '
